{"js": "// 1) Merge \"AddTransient\" + \" method.\" into a single run \"AddTransient method.\"\n//    (text itself does not change, only the run is normalized/merged).\nconst addTransientResults = context.document.body.search(\"AddTransient method.\", { matchCase: true });\naddTransientResults.load(\"text\");\nawait context.sync();\n\nif (addTransientResults.items.length > 0) {\n  addTransientResults.items[0].insertText(\"AddTransient method.\", Word.InsertLocation.replace);\n}\n\n// 2) Fix the \"List.chtml\" typo to \"List.cshtml\".\nconst typoResults = context.document.body.search(\"List.chtml\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"List.cshtml\", Word.InsertLocation.replace);\n}\n\n// 3) Append three new bullet points after the paragraph that now ends with\n//    \"...using the ItemController.\" (same bullet/list paragraph that was edited above).\nconst anchorResults = context.document.body.search(\"using the ItemController.\", { matchCase: true });\nanchorResults.load(\"text\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n\n  const newBullet1 = anchorParagraph.insertParagraph(\n    \"Added the relevant model link to the _ViewImports.cs file to import the ViewModels folder when rendering List.cshtml view template and deleted unnecessary imports from the List view template.\",\n    Word.InsertLocation.after\n  );\n\n  const newBullet2 = newBullet1.insertParagraph(\n    \"Added the bower file. Edited the MockItemRepository.cs file and List.cshtml file with correct configurations.\",\n    Word.InsertLocation.after\n  );\n\n  const newBullet3 = newBullet2.insertParagraph(\n    \"Added Bootstrap to List.cshtml and edited the _Layout.cshtml file accordingly.\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Merge \"AddTransient\" + \" method.\" into a single run \"AddTransient method.\"\n#    (rendered text is unchanged, only the run split is normalized/collapsed).\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"AddTransient method.\"\n$found = $find.Execute()\nif ($found) {\n    $lang = $range.LanguageID\n    $range.Delete()\n    $range.InsertAfter(\"AddTransient method.\")\n    $range.LanguageID = $lang\n}\n\n# 2) Fix the \"List.chtml\" typo to \"List.cshtml\".\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"List.chtml\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $lang2 = $range2.LanguageID\n    $range2.Delete()\n    $range2.InsertAfter(\"List.cshtml\")\n    $range2.LanguageID = $lang2\n}\n\n# 3) Append three new bullet points after the paragraph that now ends with\n#    \"...using the ItemController.\" (the same bullet/list paragraph edited above).\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = \"using the ItemController.\"\n$found3 = $find3.Execute()\nif ($found3) {\n    # Locate the 1-based index of the paragraph that contains the found range.\n    $anchorIndex = 0\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -le $range3.Start -and $p.Range.End -ge $range3.End) {\n            $anchorIndex = $i\n            break\n        }\n    }\n\n    if ($anchorIndex -gt 0) {\n        $anchorPara = $d.Paragraphs.Item($anchorIndex)\n        $anchorPara.Range.InsertParagraphAfter() | Out-Null\n        $newPara1 = $d.Paragraphs.Item($anchorIndex + 1)\n        $newPara1.Range.Text = \"Added the relevant model link to the _ViewImports.cs file to import the ViewModels folder when rendering List.cshtml view template and deleted unnecessary imports from the List view template.\"\n\n        $newPara1.Range.InsertParagraphAfter() | Out-Null\n        $newPara2 = $d.Paragraphs.Item($anchorIndex + 2)\n        $newPara2.Range.Text = \"Added the bower file. Edited the MockItemRepository.cs file and List.cshtml file with correct configurations.\"\n\n        $newPara2.Range.InsertParagraphAfter() | Out-Null\n        $newPara3 = $d.Paragraphs.Item($anchorIndex + 3)\n        $newPara3.Range.Text = \"Added Bootstrap to List.cshtml and edited the _Layout.cshtml file accordingly.\"\n    }\n}\n"}
